$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B values were re-keyed (primary key of the "file" table changed) ---
# Leading apostrophe forces text interpretation (keeps the existing
# quote-prefixed text style) instead of Excel auto-parsing these
# dd/mm/yyyy-looking strings as real dates.
$ws.Cells.Item(1,2).Value  = "'01/08/2023"
$ws.Cells.Item(2,2).Value  = "'01/08/2023"
$ws.Cells.Item(3,2).Value  = "'01/08/2023"
$ws.Cells.Item(4,2).Value  = "'01/08/2023"
$ws.Cells.Item(5,2).Value  = "'01/08/2023"
$ws.Cells.Item(6,2).Value  = "'01/08/2023"
$ws.Cells.Item(7,2).Value  = "'01/08/2023"
$ws.Cells.Item(8,2).Value  = "'01/08/2023"
$ws.Cells.Item(9,2).Value  = "'17/06/2023"
$ws.Cells.Item(10,2).Value = "'01/08/2023"
$ws.Cells.Item(11,2).Value = "'01/08/2023"
$ws.Cells.Item(12,2).Value = "'2023/08/01"
$ws.Cells.Item(13,2).Value = "'18/06/2023"
$ws.Cells.Item(14,2).Value = "'08/01/2023"
$ws.Cells.Item(15,2).Value = "17/06/2023"
$ws.Cells.Item(16,2).Value = "'01/08/2023"
$ws.Cells.Item(17,2).Value = "'01/08/2023"

# --- Column widths widened for C and E to fit the new values ---
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 17

# --- Selection collapsed from B1:B11 down to just B1 ---
$ws.Range("B1").Select()
